# Add a new "2022" data column (column S) to the SDG 3.4.2 (suicide
# mortality rate) sheet, mirroring the existing "2021" column (R) for
# layout/formatting, then fill in the 2022 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column R (the most recent year column, with all of its row-by-row
# formatting) and insert it as a new column immediately to its right.
# This both shifts nothing else and gives the new column S the same
# per-row styles (fonts, number formats, alignment) that column R uses.
$ws.Range("R1:R15").Copy()
$ws.Range("S1:S15").Insert(-4161)

# Header row: the new column is year 2022.
$ws.Range("S4").Value = 2022

# Data rows: suicide mortality rate per 100,000 population, 2022 figures.
$ws.Range("S5").Value = 4.9000000000000004
$ws.Range("S6").Value = 3.4
$ws.Range("S7").Value = 3.5
$ws.Range("S8").Value = 13.1
$ws.Range("S9").Value = 8.1
$ws.Range("S10").Value = 2.5
$ws.Range("S11").Value = 2.6
$ws.Range("S12").Value = 10.8
$ws.Range("S13").Value = 2.1
$ws.Range("S14").Value = 1.1000000000000001

# Match the saved selection cursor position recorded in the workbook.
$null = $ws.Range("T4").Select()
